$wb = $excel.ActiveWorkbook

# --- Dashboard sheet: latest-analysis + summary-statistics refresh ---
$ws = $wb.Worksheets.Item("Dashboard")
$ws.Range("B4").Value = "10:20:26 22/01/2026"
$ws.Range("B5").Value = "`$130.49"
$ws.Range("B6").Value = "2.05%"
$ws.Range("B7").Value = "`$130.23"
$ws.Range("B8").Value = "bullish"
$ws.Range("B9").Value = "0.049"
$ws.Range("B10").Value = "0.62"
$ws.Range("B11").Value = "Consider long position"
$ws.Range("B15").Value = 33
$ws.Range("B16").Value = "`$128.41"
$ws.Range("B18").Value = "14.84%"
$ws.Range("B19").Value = "41.56"
$ws.Range("B20").Value = "0.023"

# --- Data sheet: 7 new report rows (28-34), all 40 columns ---
$ws = $wb.Worksheets.Item("Data")
$ws.Cells.Item(28, 1).Value = "01:15:29 22/01/2026"
$ws.Cells.Item(28, 2).Value = 46044.05241941915
$ws.Cells.Item(28, 3).Value = 131.34
$ws.Cells.Item(28, 4).Value = 131.45
$ws.Cells.Item(28, 5).Value = 124.68
$ws.Cells.Item(28, 6).Value = 4.34
$ws.Cells.Item(28, 7).Value = 3.417
$ws.Cells.Item(28, 8).Value = 3129151.403
$ws.Cells.Item(28, 9).Value = 400041652.26667
$ws.Cells.Item(28, 10).Value = 1522335
$ws.Cells.Item(28, 11).Value = 4.5
$ws.Cells.Item(28, 12).Value = 5.5
$ws.Cells.Item(28, 13).Value = 4.5
$ws.Cells.Item(28, 14).Value = 7
$ws.Cells.Item(28, 15).Value = 2
$ws.Cells.Item(28, 16).Value = "neutral"
$ws.Cells.Item(28, 17).Value = 0.5
$ws.Cells.Item(28, 18).Value = 0.16
$ws.Cells.Item(28, 19).Value = "low"
$ws.Cells.Item(28, 20).Value = 58.07
$ws.Cells.Item(28, 21).Value = "neutral"
$ws.Cells.Item(28, 22).Value = -0.1001
$ws.Cells.Item(28, 23).Value = "bullish"
$ws.Cells.Item(28, 24).Value = 0.022
$ws.Cells.Item(28, 25).Value = 0.036
$ws.Cells.Item(28, 26).Value = 0
$ws.Cells.Item(28, 27).Value = 61.9
$ws.Cells.Item(28, 28).Value = "Greed"
$ws.Cells.Item(28, 29).Value = 0.552
$ws.Cells.Item(28, 30).Value = 0.448
$ws.Cells.Item(28, 31).Value = 24.05
$ws.Cells.Item(28, 32).Value = 127.91
$ws.Cells.Item(28, 33).Value = -2.61
$ws.Cells.Item(28, 34).Value = "downward"
$ws.Cells.Item(28, 35).Value = 130.22
$ws.Cells.Item(28, 36).Value = 132.28
$ws.Cells.Item(28, 37).Value = 130.22
$ws.Cells.Item(28, 38).Value = 0
$ws.Cells.Item(28, 39).Value = "Hold or wait for clearer signals"
$ws.Cells.Item(28, 40).Value = 0.55
$ws.Cells.Item(28, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(29, 1).Value = "01:20:28 22/01/2026"
$ws.Cells.Item(29, 2).Value = 46044.05588241426
$ws.Cells.Item(29, 3).Value = 131.38
$ws.Cells.Item(29, 4).Value = 131.61
$ws.Cells.Item(29, 5).Value = 124.68
$ws.Cells.Item(29, 6).Value = 4.43
$ws.Cells.Item(29, 7).Value = 3.49
$ws.Cells.Item(29, 8).Value = 3163777.966
$ws.Cells.Item(29, 9).Value = 404605773.32788
$ws.Cells.Item(29, 10).Value = 1538906
$ws.Cells.Item(29, 11).Value = 4.5
$ws.Cells.Item(29, 12).Value = 5.5
$ws.Cells.Item(29, 13).Value = 4.5
$ws.Cells.Item(29, 14).Value = 7
$ws.Cells.Item(29, 15).Value = 2
$ws.Cells.Item(29, 16).Value = "neutral"
$ws.Cells.Item(29, 17).Value = 0.5
$ws.Cells.Item(29, 18).Value = 0.16
$ws.Cells.Item(29, 19).Value = "low"
$ws.Cells.Item(29, 20).Value = 58.14
$ws.Cells.Item(29, 21).Value = "neutral"
$ws.Cells.Item(29, 22).Value = -0.0977
$ws.Cells.Item(29, 23).Value = "bullish"
$ws.Cells.Item(29, 24).Value = 0.022
$ws.Cells.Item(29, 25).Value = 0.036
$ws.Cells.Item(29, 26).Value = 0
$ws.Cells.Item(29, 27).Value = 62
$ws.Cells.Item(29, 28).Value = "Greed"
$ws.Cells.Item(29, 29).Value = 0.552
$ws.Cells.Item(29, 30).Value = 0.448
$ws.Cells.Item(29, 31).Value = 24.1
$ws.Cells.Item(29, 32).Value = 127.94
$ws.Cells.Item(29, 33).Value = -2.61
$ws.Cells.Item(29, 34).Value = "downward"
$ws.Cells.Item(29, 35).Value = 130.22
$ws.Cells.Item(29, 36).Value = 132.28
$ws.Cells.Item(29, 37).Value = 130.22
$ws.Cells.Item(29, 38).Value = 0
$ws.Cells.Item(29, 39).Value = "Hold or wait for clearer signals"
$ws.Cells.Item(29, 40).Value = 0.55
$ws.Cells.Item(29, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(30, 1).Value = "01:25:28 22/01/2026"
$ws.Cells.Item(30, 2).Value = 46044.0593591213
$ws.Cells.Item(30, 3).Value = 131.42
$ws.Cells.Item(30, 4).Value = 131.61
$ws.Cells.Item(30, 5).Value = 124.68
$ws.Cells.Item(30, 6).Value = 4.1
$ws.Cells.Item(30, 7).Value = 3.22
$ws.Cells.Item(30, 8).Value = 3192315.455
$ws.Cells.Item(30, 9).Value = 408358699.34332
$ws.Cells.Item(30, 10).Value = 1551611
$ws.Cells.Item(30, 11).Value = 4.5
$ws.Cells.Item(30, 12).Value = 5.5
$ws.Cells.Item(30, 13).Value = 4.5
$ws.Cells.Item(30, 14).Value = 7
$ws.Cells.Item(30, 15).Value = 2
$ws.Cells.Item(30, 16).Value = "neutral"
$ws.Cells.Item(30, 17).Value = 0.5
$ws.Cells.Item(30, 18).Value = 0.17
$ws.Cells.Item(30, 19).Value = "low"
$ws.Cells.Item(30, 20).Value = 58.29
$ws.Cells.Item(30, 21).Value = "neutral"
$ws.Cells.Item(30, 22).Value = -0.0929
$ws.Cells.Item(30, 23).Value = "bullish"
$ws.Cells.Item(30, 24).Value = 0.022
$ws.Cells.Item(30, 25).Value = 0.036
$ws.Cells.Item(30, 26).Value = 0
$ws.Cells.Item(30, 27).Value = 61.8
$ws.Cells.Item(30, 28).Value = "Greed"
$ws.Cells.Item(30, 29).Value = 0.552
$ws.Cells.Item(30, 30).Value = 0.448
$ws.Cells.Item(30, 31).Value = 24.21
$ws.Cells.Item(30, 32).Value = 128
$ws.Cells.Item(30, 33).Value = -2.6
$ws.Cells.Item(30, 34).Value = "downward"
$ws.Cells.Item(30, 35).Value = 130.23
$ws.Cells.Item(30, 36).Value = 132.28
$ws.Cells.Item(30, 37).Value = 130.23
$ws.Cells.Item(30, 38).Value = 0
$ws.Cells.Item(30, 39).Value = "Hold or wait for clearer signals"
$ws.Cells.Item(30, 40).Value = 0.55
$ws.Cells.Item(30, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(31, 1).Value = "01:30:28 22/01/2026"
$ws.Cells.Item(31, 2).Value = 46044.06283290812
$ws.Cells.Item(31, 3).Value = 131.52
$ws.Cells.Item(31, 4).Value = 131.68
$ws.Cells.Item(31, 5).Value = 124.68
$ws.Cells.Item(31, 6).Value = 4.45
$ws.Cells.Item(31, 7).Value = 3.502
$ws.Cells.Item(31, 8).Value = 3205652.301
$ws.Cells.Item(31, 9).Value = 410129030.24927
$ws.Cells.Item(31, 10).Value = 1559785
$ws.Cells.Item(31, 11).Value = 4.5
$ws.Cells.Item(31, 12).Value = 5.5
$ws.Cells.Item(31, 13).Value = 4.5
$ws.Cells.Item(31, 14).Value = 7
$ws.Cells.Item(31, 15).Value = 2
$ws.Cells.Item(31, 16).Value = "bullish"
$ws.Cells.Item(31, 17).Value = 0.75
$ws.Cells.Item(31, 18).Value = 0.19
$ws.Cells.Item(31, 19).Value = "low"
$ws.Cells.Item(31, 20).Value = 59.47
$ws.Cells.Item(31, 21).Value = "neutral"
$ws.Cells.Item(31, 22).Value = 0.1819
$ws.Cells.Item(31, 23).Value = "bullish"
$ws.Cells.Item(31, 24).Value = 0.022
$ws.Cells.Item(31, 25).Value = 0.036
$ws.Cells.Item(31, 26).Value = 0
$ws.Cells.Item(31, 27).Value = 62
$ws.Cells.Item(31, 28).Value = "Greed"
$ws.Cells.Item(31, 29).Value = 0.627
$ws.Cells.Item(31, 30).Value = 0.373
$ws.Cells.Item(31, 31).Value = 23.96
$ws.Cells.Item(31, 32).Value = 128.51
$ws.Cells.Item(31, 33).Value = -2.29
$ws.Cells.Item(31, 34).Value = "downward"
$ws.Cells.Item(31, 35).Value = 131.15
$ws.Cells.Item(31, 36).Value = 132.67
$ws.Cells.Item(31, 37).Value = 131.15
$ws.Cells.Item(31, 38).Value = 0.33
$ws.Cells.Item(31, 39).Value = "Consider long position"
$ws.Cells.Item(31, 40).Value = 0.59
$ws.Cells.Item(31, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(32, 1).Value = "01:35:28 22/01/2026"
$ws.Cells.Item(32, 2).Value = 46044.06629963835
$ws.Cells.Item(32, 3).Value = 131.01
$ws.Cells.Item(32, 4).Value = 131.68
$ws.Cells.Item(32, 5).Value = 124.68
$ws.Cells.Item(32, 6).Value = 3.99
$ws.Cells.Item(32, 7).Value = 3.141
$ws.Cells.Item(32, 8).Value = 3232365.927
$ws.Cells.Item(32, 9).Value = 413648315.97641
$ws.Cells.Item(32, 10).Value = 1568143
$ws.Cells.Item(32, 11).Value = 4.5
$ws.Cells.Item(32, 12).Value = 5.5
$ws.Cells.Item(32, 13).Value = 4.5
$ws.Cells.Item(32, 14).Value = 7
$ws.Cells.Item(32, 15).Value = 2
$ws.Cells.Item(32, 16).Value = "bullish"
$ws.Cells.Item(32, 17).Value = 0.75
$ws.Cells.Item(32, 18).Value = 0.16
$ws.Cells.Item(32, 19).Value = "low"
$ws.Cells.Item(32, 20).Value = 57.79
$ws.Cells.Item(32, 21).Value = "neutral"
$ws.Cells.Item(32, 22).Value = 0.1412
$ws.Cells.Item(32, 23).Value = "bullish"
$ws.Cells.Item(32, 24).Value = 0.022
$ws.Cells.Item(32, 25).Value = 0.036
$ws.Cells.Item(32, 26).Value = 0
$ws.Cells.Item(32, 27).Value = 61.7
$ws.Cells.Item(32, 28).Value = "Greed"
$ws.Cells.Item(32, 29).Value = 0.627
$ws.Cells.Item(32, 30).Value = 0.373
$ws.Cells.Item(32, 31).Value = 24.07
$ws.Cells.Item(32, 32).Value = 127.97
$ws.Cells.Item(32, 33).Value = -2.32
$ws.Cells.Item(32, 34).Value = "downward"
$ws.Cells.Item(32, 35).Value = 131.01
$ws.Cells.Item(32, 36).Value = 132.67
$ws.Cells.Item(32, 37).Value = 131.01
$ws.Cells.Item(32, 38).Value = 0.33
$ws.Cells.Item(32, 39).Value = "Consider long position"
$ws.Cells.Item(32, 40).Value = 0.59
$ws.Cells.Item(32, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(33, 1).Value = "01:40:27 22/01/2026"
$ws.Cells.Item(33, 2).Value = 46044.06976695727
$ws.Cells.Item(33, 3).Value = 130.61
$ws.Cells.Item(33, 4).Value = 131.68
$ws.Cells.Item(33, 5).Value = 124.68
$ws.Cells.Item(33, 6).Value = 3.34
$ws.Cells.Item(33, 7).Value = 2.624
$ws.Cells.Item(33, 8).Value = 3243674.964
$ws.Cells.Item(33, 9).Value = 415141055.79003
$ws.Cells.Item(33, 10).Value = 1576177
$ws.Cells.Item(33, 11).Value = 4.5
$ws.Cells.Item(33, 12).Value = 5.5
$ws.Cells.Item(33, 13).Value = 4.5
$ws.Cells.Item(33, 14).Value = 7
$ws.Cells.Item(33, 15).Value = 2
$ws.Cells.Item(33, 16).Value = "bullish"
$ws.Cells.Item(33, 17).Value = 0.75
$ws.Cells.Item(33, 18).Value = 0.13
$ws.Cells.Item(33, 19).Value = "low"
$ws.Cells.Item(33, 20).Value = 56.51
$ws.Cells.Item(33, 21).Value = "neutral"
$ws.Cells.Item(33, 22).Value = 0.1093
$ws.Cells.Item(33, 23).Value = "bullish"
$ws.Cells.Item(33, 24).Value = 0.022
$ws.Cells.Item(33, 25).Value = 0.036
$ws.Cells.Item(33, 26).Value = 0
$ws.Cells.Item(33, 27).Value = 61.3
$ws.Cells.Item(33, 28).Value = "Greed"
$ws.Cells.Item(33, 29).Value = 0.627
$ws.Cells.Item(33, 30).Value = 0.373
$ws.Cells.Item(33, 31).Value = 24.23
$ws.Cells.Item(33, 32).Value = 127.55
$ws.Cells.Item(33, 33).Value = -2.34
$ws.Cells.Item(33, 34).Value = "downward"
$ws.Cells.Item(33, 35).Value = 129.57
$ws.Cells.Item(33, 36).Value = 130.91
$ws.Cells.Item(33, 37).Value = 129.57
$ws.Cells.Item(33, 38).Value = 0.33
$ws.Cells.Item(33, 39).Value = "Consider long position"
$ws.Cells.Item(33, 40).Value = 0.59
$ws.Cells.Item(33, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(34, 1).Value = "10:20:26 22/01/2026"
$ws.Cells.Item(34, 2).Value = 46044.4308623914
$ws.Cells.Item(34, 3).Value = 130.49
$ws.Cells.Item(34, 4).Value = 132.17
$ws.Cells.Item(34, 5).Value = 125.26
$ws.Cells.Item(34, 6).Value = 2.62
$ws.Cells.Item(34, 7).Value = 2.049
$ws.Cells.Item(34, 8).Value = 3076104.036
$ws.Cells.Item(34, 9).Value = 396518874.96492
$ws.Cells.Item(34, 10).Value = 1542990
$ws.Cells.Item(34, 11).Value = 4.5
$ws.Cells.Item(34, 12).Value = 5.5
$ws.Cells.Item(34, 13).Value = 4.5
$ws.Cells.Item(34, 14).Value = 7
$ws.Cells.Item(34, 15).Value = 2
$ws.Cells.Item(34, 16).Value = "bullish"
$ws.Cells.Item(34, 17).Value = 1
$ws.Cells.Item(34, 18).Value = 0.06
$ws.Cells.Item(34, 19).Value = "medium"
$ws.Cells.Item(34, 20).Value = 53.04
$ws.Cells.Item(34, 21).Value = "neutral"
$ws.Cells.Item(34, 22).Value = 0.4808
$ws.Cells.Item(34, 23).Value = "bullish"
$ws.Cells.Item(34, 24).Value = 0.049
$ws.Cells.Item(34, 25).Value = 0.082
$ws.Cells.Item(34, 26).Value = 0
$ws.Cells.Item(34, 27).Value = 52.4
$ws.Cells.Item(34, 28).Value = "Neutral"
$ws.Cells.Item(34, 29).Value = 0.655
$ws.Cells.Item(34, 30).Value = 0.345
$ws.Cells.Item(34, 31).Value = 24.13
$ws.Cells.Item(34, 32).Value = 130.23
$ws.Cells.Item(34, 33).Value = -0.2
$ws.Cells.Item(34, 34).Value = "downward"
$ws.Cells.Item(34, 35).Value = 129.24
$ws.Cells.Item(34, 36).Value = 130.65
$ws.Cells.Item(34, 37).Value = 129.24
$ws.Cells.Item(34, 38).Value = 0.33
$ws.Cells.Item(34, 39).Value = "Consider long position"
$ws.Cells.Item(34, 40).Value = 0.62
$ws.Cells.Item(34, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

# --- Price Analysis sheet: 7 new rows (30-36), mirrored from Data ---
$ws = $wb.Worksheets.Item("Price Analysis")
$ws.Cells.Item(30, 1).Value = "01:15:29 22/01/2026"
$ws.Cells.Item(30, 2).Value = 46044.05241941915
$ws.Cells.Item(30, 3).Value = 131.34
$ws.Cells.Item(30, 4).Value = 131.45
$ws.Cells.Item(30, 5).Value = 124.68
$ws.Cells.Item(30, 6).Value = 4.34
$ws.Cells.Item(30, 7).Value = 3.417
$ws.Cells.Item(30, 8).Value = 3129151.403
$ws.Cells.Item(30, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(31, 1).Value = "01:20:28 22/01/2026"
$ws.Cells.Item(31, 2).Value = 46044.05588241426
$ws.Cells.Item(31, 3).Value = 131.38
$ws.Cells.Item(31, 4).Value = 131.61
$ws.Cells.Item(31, 5).Value = 124.68
$ws.Cells.Item(31, 6).Value = 4.43
$ws.Cells.Item(31, 7).Value = 3.49
$ws.Cells.Item(31, 8).Value = 3163777.966
$ws.Cells.Item(31, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(32, 1).Value = "01:25:28 22/01/2026"
$ws.Cells.Item(32, 2).Value = 46044.0593591213
$ws.Cells.Item(32, 3).Value = 131.42
$ws.Cells.Item(32, 4).Value = 131.61
$ws.Cells.Item(32, 5).Value = 124.68
$ws.Cells.Item(32, 6).Value = 4.1
$ws.Cells.Item(32, 7).Value = 3.22
$ws.Cells.Item(32, 8).Value = 3192315.455
$ws.Cells.Item(32, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(33, 1).Value = "01:30:28 22/01/2026"
$ws.Cells.Item(33, 2).Value = 46044.06283290812
$ws.Cells.Item(33, 3).Value = 131.52
$ws.Cells.Item(33, 4).Value = 131.68
$ws.Cells.Item(33, 5).Value = 124.68
$ws.Cells.Item(33, 6).Value = 4.45
$ws.Cells.Item(33, 7).Value = 3.502
$ws.Cells.Item(33, 8).Value = 3205652.301
$ws.Cells.Item(33, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(34, 1).Value = "01:35:28 22/01/2026"
$ws.Cells.Item(34, 2).Value = 46044.06629963835
$ws.Cells.Item(34, 3).Value = 131.01
$ws.Cells.Item(34, 4).Value = 131.68
$ws.Cells.Item(34, 5).Value = 124.68
$ws.Cells.Item(34, 6).Value = 3.99
$ws.Cells.Item(34, 7).Value = 3.141
$ws.Cells.Item(34, 8).Value = 3232365.927
$ws.Cells.Item(34, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(35, 1).Value = "01:40:27 22/01/2026"
$ws.Cells.Item(35, 2).Value = 46044.06976695727
$ws.Cells.Item(35, 3).Value = 130.61
$ws.Cells.Item(35, 4).Value = 131.68
$ws.Cells.Item(35, 5).Value = 124.68
$ws.Cells.Item(35, 6).Value = 3.34
$ws.Cells.Item(35, 7).Value = 2.624
$ws.Cells.Item(35, 8).Value = 3243674.964
$ws.Cells.Item(35, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(36, 1).Value = "10:20:26 22/01/2026"
$ws.Cells.Item(36, 2).Value = 46044.4308623914
$ws.Cells.Item(36, 3).Value = 130.49
$ws.Cells.Item(36, 4).Value = 132.17
$ws.Cells.Item(36, 5).Value = 125.26
$ws.Cells.Item(36, 6).Value = 2.62
$ws.Cells.Item(36, 7).Value = 2.049
$ws.Cells.Item(36, 8).Value = 3076104.036
$ws.Cells.Item(36, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

# --- Technical Analysis sheet: 7 new rows (30-36), mirrored from Data ---
$ws = $wb.Worksheets.Item("Technical Analysis")
$ws.Cells.Item(30, 1).Value = "01:15:29 22/01/2026"
$ws.Cells.Item(30, 2).Value = 46044.05241941915
$ws.Cells.Item(30, 3).Value = 0.5
$ws.Cells.Item(30, 4).Value = 0.16
$ws.Cells.Item(30, 5).Value = 58.07
$ws.Cells.Item(30, 6).Value = -0.1001
$ws.Cells.Item(30, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(31, 1).Value = "01:20:28 22/01/2026"
$ws.Cells.Item(31, 2).Value = 46044.05588241426
$ws.Cells.Item(31, 3).Value = 0.5
$ws.Cells.Item(31, 4).Value = 0.16
$ws.Cells.Item(31, 5).Value = 58.14
$ws.Cells.Item(31, 6).Value = -0.0977
$ws.Cells.Item(31, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(32, 1).Value = "01:25:28 22/01/2026"
$ws.Cells.Item(32, 2).Value = 46044.0593591213
$ws.Cells.Item(32, 3).Value = 0.5
$ws.Cells.Item(32, 4).Value = 0.17
$ws.Cells.Item(32, 5).Value = 58.29
$ws.Cells.Item(32, 6).Value = -0.0929
$ws.Cells.Item(32, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(33, 1).Value = "01:30:28 22/01/2026"
$ws.Cells.Item(33, 2).Value = 46044.06283290812
$ws.Cells.Item(33, 3).Value = 0.75
$ws.Cells.Item(33, 4).Value = 0.19
$ws.Cells.Item(33, 5).Value = 59.47
$ws.Cells.Item(33, 6).Value = 0.1819
$ws.Cells.Item(33, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(34, 1).Value = "01:35:28 22/01/2026"
$ws.Cells.Item(34, 2).Value = 46044.06629963835
$ws.Cells.Item(34, 3).Value = 0.75
$ws.Cells.Item(34, 4).Value = 0.16
$ws.Cells.Item(34, 5).Value = 57.79
$ws.Cells.Item(34, 6).Value = 0.1412
$ws.Cells.Item(34, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(35, 1).Value = "01:40:27 22/01/2026"
$ws.Cells.Item(35, 2).Value = 46044.06976695727
$ws.Cells.Item(35, 3).Value = 0.75
$ws.Cells.Item(35, 4).Value = 0.13
$ws.Cells.Item(35, 5).Value = 56.51
$ws.Cells.Item(35, 6).Value = 0.1093
$ws.Cells.Item(35, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(36, 1).Value = "10:20:26 22/01/2026"
$ws.Cells.Item(36, 2).Value = 46044.4308623914
$ws.Cells.Item(36, 3).Value = 1
$ws.Cells.Item(36, 4).Value = 0.06
$ws.Cells.Item(36, 5).Value = 53.04
$ws.Cells.Item(36, 6).Value = 0.4808
$ws.Cells.Item(36, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

# --- Fundamental Analysis sheet: 7 new rows (30-36), mirrored from Data ---
$ws = $wb.Worksheets.Item("Fundamental Analysis")
$ws.Cells.Item(30, 1).Value = "01:15:29 22/01/2026"
$ws.Cells.Item(30, 2).Value = 46044.05241941915
$ws.Cells.Item(30, 3).Value = 4.5
$ws.Cells.Item(30, 4).Value = 5.5
$ws.Cells.Item(30, 5).Value = 4.5
$ws.Cells.Item(30, 6).Value = 7
$ws.Cells.Item(30, 7).Value = 2
$ws.Cells.Item(30, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(31, 1).Value = "01:20:28 22/01/2026"
$ws.Cells.Item(31, 2).Value = 46044.05588241426
$ws.Cells.Item(31, 3).Value = 4.5
$ws.Cells.Item(31, 4).Value = 5.5
$ws.Cells.Item(31, 5).Value = 4.5
$ws.Cells.Item(31, 6).Value = 7
$ws.Cells.Item(31, 7).Value = 2
$ws.Cells.Item(31, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(32, 1).Value = "01:25:28 22/01/2026"
$ws.Cells.Item(32, 2).Value = 46044.0593591213
$ws.Cells.Item(32, 3).Value = 4.5
$ws.Cells.Item(32, 4).Value = 5.5
$ws.Cells.Item(32, 5).Value = 4.5
$ws.Cells.Item(32, 6).Value = 7
$ws.Cells.Item(32, 7).Value = 2
$ws.Cells.Item(32, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(33, 1).Value = "01:30:28 22/01/2026"
$ws.Cells.Item(33, 2).Value = 46044.06283290812
$ws.Cells.Item(33, 3).Value = 4.5
$ws.Cells.Item(33, 4).Value = 5.5
$ws.Cells.Item(33, 5).Value = 4.5
$ws.Cells.Item(33, 6).Value = 7
$ws.Cells.Item(33, 7).Value = 2
$ws.Cells.Item(33, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(34, 1).Value = "01:35:28 22/01/2026"
$ws.Cells.Item(34, 2).Value = 46044.06629963835
$ws.Cells.Item(34, 3).Value = 4.5
$ws.Cells.Item(34, 4).Value = 5.5
$ws.Cells.Item(34, 5).Value = 4.5
$ws.Cells.Item(34, 6).Value = 7
$ws.Cells.Item(34, 7).Value = 2
$ws.Cells.Item(34, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(35, 1).Value = "01:40:27 22/01/2026"
$ws.Cells.Item(35, 2).Value = 46044.06976695727
$ws.Cells.Item(35, 3).Value = 4.5
$ws.Cells.Item(35, 4).Value = 5.5
$ws.Cells.Item(35, 5).Value = 4.5
$ws.Cells.Item(35, 6).Value = 7
$ws.Cells.Item(35, 7).Value = 2
$ws.Cells.Item(35, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(36, 1).Value = "10:20:26 22/01/2026"
$ws.Cells.Item(36, 2).Value = 46044.4308623914
$ws.Cells.Item(36, 3).Value = 4.5
$ws.Cells.Item(36, 4).Value = 5.5
$ws.Cells.Item(36, 5).Value = 4.5
$ws.Cells.Item(36, 6).Value = 7
$ws.Cells.Item(36, 7).Value = 2
$ws.Cells.Item(36, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

# --- Sentiment Analysis sheet: 7 new rows (30-36), mirrored from Data ---
$ws = $wb.Worksheets.Item("Sentiment Analysis")
$ws.Cells.Item(30, 1).Value = "01:15:29 22/01/2026"
$ws.Cells.Item(30, 2).Value = 46044.05241941915
$ws.Cells.Item(30, 3).Value = 0.022
$ws.Cells.Item(30, 4).Value = 0.036
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = 61.9
$ws.Cells.Item(30, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(31, 1).Value = "01:20:28 22/01/2026"
$ws.Cells.Item(31, 2).Value = 46044.05588241426
$ws.Cells.Item(31, 3).Value = 0.022
$ws.Cells.Item(31, 4).Value = 0.036
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = 62
$ws.Cells.Item(31, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(32, 1).Value = "01:25:28 22/01/2026"
$ws.Cells.Item(32, 2).Value = 46044.0593591213
$ws.Cells.Item(32, 3).Value = 0.022
$ws.Cells.Item(32, 4).Value = 0.036
$ws.Cells.Item(32, 5).Value = 0
$ws.Cells.Item(32, 6).Value = 61.8
$ws.Cells.Item(32, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(33, 1).Value = "01:30:28 22/01/2026"
$ws.Cells.Item(33, 2).Value = 46044.06283290812
$ws.Cells.Item(33, 3).Value = 0.022
$ws.Cells.Item(33, 4).Value = 0.036
$ws.Cells.Item(33, 5).Value = 0
$ws.Cells.Item(33, 6).Value = 62
$ws.Cells.Item(33, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(34, 1).Value = "01:35:28 22/01/2026"
$ws.Cells.Item(34, 2).Value = 46044.06629963835
$ws.Cells.Item(34, 3).Value = 0.022
$ws.Cells.Item(34, 4).Value = 0.036
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 6).Value = 61.7
$ws.Cells.Item(34, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(35, 1).Value = "01:40:27 22/01/2026"
$ws.Cells.Item(35, 2).Value = 46044.06976695727
$ws.Cells.Item(35, 3).Value = 0.022
$ws.Cells.Item(35, 4).Value = 0.036
$ws.Cells.Item(35, 5).Value = 0
$ws.Cells.Item(35, 6).Value = 61.3
$ws.Cells.Item(35, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(36, 1).Value = "10:20:26 22/01/2026"
$ws.Cells.Item(36, 2).Value = 46044.4308623914
$ws.Cells.Item(36, 3).Value = 0.049
$ws.Cells.Item(36, 4).Value = 0.082
$ws.Cells.Item(36, 5).Value = 0
$ws.Cells.Item(36, 6).Value = 52.4
$ws.Cells.Item(36, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

# --- Predictions sheet: 7 new rows (30-36), mirrored from Data ---
$ws = $wb.Worksheets.Item("Predictions")
$ws.Cells.Item(30, 1).Value = "01:15:29 22/01/2026"
$ws.Cells.Item(30, 2).Value = 46044.05241941915
$ws.Cells.Item(30, 3).Value = 131.34
$ws.Cells.Item(30, 4).Value = 127.91
$ws.Cells.Item(30, 5).Value = -2.61
$ws.Cells.Item(30, 6).Value = 0.552
$ws.Cells.Item(30, 7).Value = 0.448
$ws.Cells.Item(30, 8).Value = 130.22
$ws.Cells.Item(30, 9).Value = 132.28
$ws.Cells.Item(30, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(31, 1).Value = "01:20:28 22/01/2026"
$ws.Cells.Item(31, 2).Value = 46044.05588241426
$ws.Cells.Item(31, 3).Value = 131.38
$ws.Cells.Item(31, 4).Value = 127.94
$ws.Cells.Item(31, 5).Value = -2.61
$ws.Cells.Item(31, 6).Value = 0.552
$ws.Cells.Item(31, 7).Value = 0.448
$ws.Cells.Item(31, 8).Value = 130.22
$ws.Cells.Item(31, 9).Value = 132.28
$ws.Cells.Item(31, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(32, 1).Value = "01:25:28 22/01/2026"
$ws.Cells.Item(32, 2).Value = 46044.0593591213
$ws.Cells.Item(32, 3).Value = 131.42
$ws.Cells.Item(32, 4).Value = 128
$ws.Cells.Item(32, 5).Value = -2.6
$ws.Cells.Item(32, 6).Value = 0.552
$ws.Cells.Item(32, 7).Value = 0.448
$ws.Cells.Item(32, 8).Value = 130.23
$ws.Cells.Item(32, 9).Value = 132.28
$ws.Cells.Item(32, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(33, 1).Value = "01:30:28 22/01/2026"
$ws.Cells.Item(33, 2).Value = 46044.06283290812
$ws.Cells.Item(33, 3).Value = 131.52
$ws.Cells.Item(33, 4).Value = 128.51
$ws.Cells.Item(33, 5).Value = -2.29
$ws.Cells.Item(33, 6).Value = 0.627
$ws.Cells.Item(33, 7).Value = 0.373
$ws.Cells.Item(33, 8).Value = 131.15
$ws.Cells.Item(33, 9).Value = 132.67
$ws.Cells.Item(33, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(34, 1).Value = "01:35:28 22/01/2026"
$ws.Cells.Item(34, 2).Value = 46044.06629963835
$ws.Cells.Item(34, 3).Value = 131.01
$ws.Cells.Item(34, 4).Value = 127.97
$ws.Cells.Item(34, 5).Value = -2.32
$ws.Cells.Item(34, 6).Value = 0.627
$ws.Cells.Item(34, 7).Value = 0.373
$ws.Cells.Item(34, 8).Value = 131.01
$ws.Cells.Item(34, 9).Value = 132.67
$ws.Cells.Item(34, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(35, 1).Value = "01:40:27 22/01/2026"
$ws.Cells.Item(35, 2).Value = 46044.06976695727
$ws.Cells.Item(35, 3).Value = 130.61
$ws.Cells.Item(35, 4).Value = 127.55
$ws.Cells.Item(35, 5).Value = -2.34
$ws.Cells.Item(35, 6).Value = 0.627
$ws.Cells.Item(35, 7).Value = 0.373
$ws.Cells.Item(35, 8).Value = 129.57
$ws.Cells.Item(35, 9).Value = 130.91
$ws.Cells.Item(35, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Cells.Item(36, 1).Value = "10:20:26 22/01/2026"
$ws.Cells.Item(36, 2).Value = 46044.4308623914
$ws.Cells.Item(36, 3).Value = 130.49
$ws.Cells.Item(36, 4).Value = 130.23
$ws.Cells.Item(36, 5).Value = -0.2
$ws.Cells.Item(36, 6).Value = 0.655
$ws.Cells.Item(36, 7).Value = 0.345
$ws.Cells.Item(36, 8).Value = 129.24
$ws.Cells.Item(36, 9).Value = 130.65
$ws.Cells.Item(36, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

# --- Charts: extend every series range from row 29 to row 36 ---
foreach ($sheetName in @("Price Analysis","Technical Analysis","Fundamental Analysis","Sentiment Analysis","Predictions")) {
    $cws = $wb.Worksheets.Item($sheetName)
    $n = $cws.ChartObjects().Count
    for ($i = 1; $i -le $n; $i++) {
        $chart = $cws.ChartObjects($i).Chart
        $sc = $chart.SeriesCollection().Count
        for ($j = 1; $j -le $sc; $j++) {
            $ser = $chart.SeriesCollection($j)
            $ser.Formula = $ser.Formula.Replace("`$29", "`$36")
        }
    }
}

Write-Host "edit.ps1 complete"
